$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.69000000000058
$ws.Range("H2").Value = [double]"1.323663814754285e-16"
$ws.Range("K2").Value = 47.72475545989161
$ws.Range("L2").Value = "[41.28179077641309, 54.16772014337013]"
$ws.Range("O2").Value = 1.515763422452733
$ws.Range("P2").Value = "[1.3773949772495788, 1.6541318676558872]"
$ws.Range("S2").Value = 53.91175305302956
$ws.Range("T2").Value = "[49.651680908789416, 58.1718251972697]"
$ws.Range("W2").Value = 19.49251251251295
$ws.Range("X2").Value = 18.92676676676719
$ws.Range("Y2").Value = 20.05825825825871

# Row 3
$ws.Range("E3").Value = 23.3100000000002
$ws.Range("H3").Value = [double]"1.323663814754285e-16"
$ws.Range("K3").Value = 48.11414621162205
$ws.Range("L3").Value = "[38.782723059292806, 57.44556936395129]"
$ws.Range("O3").Value = -1.044052813805617
$ws.Range("P3").Value = "[-1.2327370572644636, -0.855368570346771]"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 53.06333937020077
$ws.Range("T3").Value = "[47.99511999513807, 58.13155874526346]"
$ws.Range("W3").Value = 3.873333333333367
$ws.Range("X3").Value = 3.173333333333361
$ws.Range("Y3").Value = 4.573333333333372
